# Update R and Datasets
# Applies the edits made to "R Scripts/Dataset/CRQ2_1_UNIQUE.xlsx":
#   - Rename column headers (B1, C1) to the new standardized names
#   - Move the active cell selection from A2 to B3
#   - Widen columns B and C to fit the new, longer header text
#   - Update the saved window position/size

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Rename the header cells with the updated text
$ws.Range("B1").Value = "NUM_MERGES_PER_100_COMMITS_WO_FT"
$ws.Range("C1").Value = "NUM_MERGES_PER_100_COMMITS_WITH_FT"

# Resize columns B and C to fit the new header text
$ws.Columns.Item(2).ColumnWidth = 42.16666666666667
$ws.Columns.Item(3).ColumnWidth = 43.66666666666667

# Move the active selection to B3
$ws.Range("B3").Select()

# Update the workbook window size/position to match the saved state
$win = $wb.Windows.Item(1)
$win.Left = 1530
$win.Top = 2130
$win.Width = 15375
$win.Height = 7875

$wb.Save()
